# "Update with work from laptop"
#
# The Bones sheet gains a new bone "Radius_Link" inserted (in terms of
# displayed order) between "Humerus" and "Radius". The Joints sheet gains a
# new joint "Wrist_Flexion_Rotator" and the old single "Elbow"/"Wrist" rows
# are replaced by three rows: "Elbow_Flexion_Rotator", "Elbow_Supination_
# Rotator" and "Wrist_Flexion_Rotator". In both sheets the literal ID (Bones
# column A) / formula (Joints column A) sequences are untouched per-row; only
# the names in column B shift down to make room, with the final, displaced
# entry landing in a freshly appended last row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Bones": shift names in B6:B32 down by one row to make room for the
# new "Radius_Link" bone, then append the row that falls off the end.
# ---------------------------------------------------------------------
$bones = $wb.Worksheets.Item("Bones")

$boneNames = @(
  "Shoulder_Horizontal_Link","Shoulder_Frontal_Link","Shoulder_Sagittal_Link","Humerus","Radius","Carpus",
  "CMC_I_Link","Metacarpus_I","MCP_I_Link","Proximal_Phalanx_I","Distal_Phalanx_I",
  "CMC_II_Link","Metacarpus_II","Proximal_Phalanx_II","Medial_Phalanx_II","Distal_Phalanx_II",
  "CMC_III_Link","Metacarpus_III","Proximal_Phalanx_III","Medial_Phalanx_III","Distal_Phalanx_III",
  "CMC_IV_Link","Metacarpus_IV","Proximal_Phalanx_IV","Medial_Phalanx_IV","Distal_Phalanx_IV",
  "CMC_V_Link","Metacarpus_V","Proximal_Phalanx","Medial_Phalanx_V","Distal_Phalanx_V"
)

# Row r (2..32) currently holds $boneNames[r-2]. Walk from the bottom up so we
# never clobber a value before it has been read.
for ($r = 32; $r -ge 7; $r--) {
  $bones.Range("B$r").Value = $boneNames[$r - 3]
}
$bones.Range("B6").Value = "Radius_Link"

# The row that fell off the bottom ("Distal_Phalanx_V") becomes a new row 33,
# continuing the same ID progression (+2 each row) and Base=FALSE pattern.
$bones.Range("A33").Value = 63
$bones.Range("B33").Value = "Distal_Phalanx_V"
$bones.Range("C33").Value = $false

# ---------------------------------------------------------------------
# Sheet "Joints": same kind of shift. Rows 5:31 shift their names down by
# one to make room for the Elbow/Wrist rows being split into three rows,
# and the row that falls off the end is appended as a new row 32 (with its
# usual "=Bones!A{row}+1" formula in column A).
# ---------------------------------------------------------------------
$joints = $wb.Worksheets.Item("Joints")

$jointNames = @(
  "Shoulder_Horizontal_Rotator","Shoulder_Frontal_Rotator","Shoulder_Sagittal_Rotator",
  "Elbow","Wrist",
  "CMC_I_Wiggle","CMC_I_Rotator","MCP_I_Swivel","MCP_I_Rotator","DIP_I_Rotator",
  "CMC_II_Wiggle","CMC_II_Rotator","MCP_II_Rotator","PIP_II_Rotator","DIP_II_Rotator",
  "CMC_III_Wiggle","CMC_III_Rotator","MCP_III_Rotator","PIP_III_Rotator","DIP_III_Rotator",
  "CMC_IV_Wiggle","CMC_IV_Rotator","MCP_IV_Rotator","PIP_IV_Rotator","DIP_IV_Rotator",
  "CMC_V_Wiggle","CMC_V_Rotator","MCP_V_Rotator","PIP_V_Rotator","DIP_V_Rotator"
)

for ($r = 31; $r -ge 8; $r--) {
  $joints.Range("B$r").Value = $jointNames[$r - 3]
}
# Write in the same order the new unique strings were first introduced
# (Elbow_Flexion_Rotator, then Wrist_Flexion_Rotator, then Elbow_Supination_
# Rotator) so the shared-string table indices line up with the target file.
$joints.Range("B5").Value = "Elbow_Flexion_Rotator"
$joints.Range("B7").Value = "Wrist_Flexion_Rotator"
$joints.Range("B6").Value = "Elbow_Supination_Rotator"

$joints.Range("A32").Formula = "=Bones!A32+1"
$joints.Range("B32").Value = "DIP_V_Rotator"

# ---------------------------------------------------------------------
# Selections left by the editor in the saved file.
# ---------------------------------------------------------------------
$bones.Activate() | Out-Null
$bones.Range("F27").Select() | Out-Null

$joints.Activate() | Out-Null
$joints.Range("B7").Select() | Out-Null
